$wb = $excel.ActiveWorkbook

$stacked = $wb.Worksheets.Item("stacked")
$stacked.Range("A1").Value = "SampleID"
$stacked.Range("B1").Value = "GCRunID"
$stacked.Range("E1").Value = "ProcessingMethod"
$stacked.Activate() | Out-Null
$stacked.Range("A2").Select() | Out-Null

$tidy = $wb.Worksheets.Item("tidy")
$tidy.Range("A1").Value = "SampleID"
$tidy.Range("B1").Value = "GCRunID"
$tidy.Range("E1").Value = "ProcessingMethod"
$tidy.Activate() | Out-Null
$tidy.Range("A2").Select() | Out-Null
